$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.914.72"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.916.38"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'199.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.44%  "
$ws.Range("D6").Value = "'599.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D8").Value = "'0.549"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").Value = "2.914.37"
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("D11").Value = "'0.440"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +18.88%  "
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'4.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "3.453.19"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "75.791.80"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "'0.0000190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "'27.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").Value = "2.913.28"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("D19").Value = "'8.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'12.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.62%  "
$ws.Range("D21").Value = "'378.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").Value = "'4.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "'71.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "3.065.81"
$ws.Range("E26").Value = "  +3.33%  "
$ws.Range("D27").Value = "'4.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "'9.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D31").Value = "'1.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").Value = "'506.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "'7.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'164.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").Value = "'20.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").Value = "'0.109"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +27.29%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "'19.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "'0.114"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.23%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'181.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "'0.344"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.54%  "
$ws.Range("D44").Value = "'5.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'40.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Value = "'0.572"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'0.663"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("D51").Value = "'3.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
